# repull data, push all data, mean calculation
# Update the dSF column (F) values for the rows whose underlying data was
# re-pulled / recalculated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 0
    7  = 1
    11 = -1
    12 = 1
    14 = 1
    26 = 2
    32 = 0
    36 = 1
    37 = 0
    40 = -1
    56 = -9
    58 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
